$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GV_long")
$ws.Activate()
for ($r = 128; $r -le 148; $r++) {
    $ws.Cells.Item($r, 3).Value = "NTU"
}
$ws.Range("C:C").Select()
